# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$ws1.Range("F2").Value = "2021-10-05 14:22:55.230941"
$ws1.Range("F3").Value = "2021-10-05 14:22:55.230948"
$ws1.Range("F4").Value = "2021-10-05 14:22:55.230952"
$ws1.Range("F5").Value = "2021-10-05 14:22:55.230954"
$ws1.Range("F6").Value = "2021-10-05 14:22:55.230957"
$ws1.Range("F7").Value = "2021-10-05 14:22:55.230959"
$ws1.Range("F8").Value = "2021-10-05 14:22:55.230962"
$ws1.Range("F9").Value = "2021-10-05 14:22:55.230964"
$ws1.Range("F10").Value = "2021-10-05 14:22:55.230967"
$ws1.Range("F11").Value = "2021-10-05 14:22:55.230970"
$ws1.Range("F12").Value = "2021-10-05 14:22:55.230972"
$ws1.Range("F13").Value = "2021-10-05 14:22:55.230975"
$ws1.Range("F14").Value = "2021-10-05 14:22:55.230977"
$ws1.Range("F15").Value = "2021-10-05 14:22:55.230979"
$ws1.Range("F16").Value = "2021-10-05 14:22:55.230982"
$ws1.Range("F17").Value = "2021-10-05 14:22:55.230984"
$ws1.Range("F18").Value = "2021-10-05 14:22:55.230987"
$ws1.Range("F19").Value = "2021-10-05 14:22:55.230990"
$ws1.Range("F20").Value = "2021-10-05 14:22:55.230992"
$ws1.Range("F21").Value = "2021-10-05 14:22:55.230995"
$ws1.Range("F22").Value = "2021-10-05 14:22:55.230997"
$ws1.Range("F23").Value = "2021-10-05 14:22:55.230999"
$ws1.Range("F24").Value = "2021-10-05 14:22:55.231002"
$ws1.Range("F25").Value = "2021-10-05 14:22:55.231004"
$ws1.Range("F26").Value = "2021-10-05 14:22:55.231007"
$ws1.Range("F27").Value = "2021-10-05 14:22:55.231010"
$ws1.Range("F28").Value = "2021-10-05 14:22:55.231012"
$ws1.Range("F29").Value = "2021-10-05 14:22:55.231014"
$ws1.Range("F30").Value = "2021-10-05 14:22:55.231017"
$ws1.Range("F31").Value = "2021-10-05 14:22:55.231019"
$ws1.Range("F32").Value = "2021-10-05 14:22:55.231022"
$ws1.Range("F33").Value = "2021-10-05 14:22:55.231024"
$ws1.Range("F34").Value = "2021-10-05 14:22:55.231027"
$ws1.Range("F35").Value = "2021-10-05 14:22:55.231029"
$ws1.Range("F36").Value = "2021-10-05 14:22:55.231032"
$ws1.Range("F37").Value = "2021-10-05 14:22:55.231034"
$ws1.Range("F38").Value = "2021-10-05 14:22:55.231037"
$ws1.Range("F39").Value = "2021-10-05 14:22:55.231039"
$ws1.Range("F40").Value = "2021-10-05 14:22:55.231041"
$ws1.Range("F41").Value = "2021-10-05 14:22:55.231044"
$ws1.Range("F42").Value = "2021-10-05 14:22:55.231047"
$ws1.Range("F43").Value = "2021-10-05 14:22:55.231049"
$ws1.Range("F44").Value = "2021-10-05 14:22:55.231052"
$ws1.Range("F45").Value = "2021-10-05 14:22:55.231054"
$ws1.Range("F46").Value = "2021-10-05 14:22:55.231056"
$ws1.Range("F47").Value = "2021-10-05 14:22:55.231059"
$ws1.Range("F48").Value = "2021-10-05 14:22:55.231061"
$ws1.Range("F49").Value = "2021-10-05 14:22:55.231064"
$ws1.Range("F50").Value = "2021-10-05 14:22:55.231066"
$ws1.Range("F51").Value = "2021-10-05 14:22:55.231068"
$ws1.Range("F52").Value = "2021-10-05 14:22:55.231071"
$ws1.Range("F53").Value = "2021-10-05 14:22:55.231073"
$ws1.Range("F54").Value = "2021-10-05 14:22:55.231076"
$ws1.Range("F55").Value = "2021-10-05 14:22:55.231078"
$ws1.Range("F56").Value = "2021-10-05 14:22:55.231081"
$ws1.Range("F57").Value = "2021-10-05 14:22:55.231083"
$ws1.Range("F58").Value = "2021-10-05 14:22:55.231086"
$ws1.Range("F59").Value = "2021-10-05 14:22:55.231088"
$ws1.Range("F60").Value = "2021-10-05 14:22:55.231090"
$ws1.Range("F61").Value = "2021-10-05 14:22:55.231093"
$ws1.Range("F62").Value = "2021-10-05 14:22:55.231095"
$ws1.Range("F63").Value = "2021-10-05 14:22:55.231097"
$ws1.Range("F64").Value = "2021-10-05 14:22:55.231100"
$ws1.Range("F65").Value = "2021-10-05 14:22:55.231102"
$ws1.Range("F66").Value = "2021-10-05 14:22:55.231106"
$ws1.Range("F67").Value = "2021-10-05 14:22:55.231108"

# --- Add the new "metadata" sheet after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Thoracic aortic aneurysm and dissection"
$ws2.Range("C2").Value = 700
$ws2.Range("E2").Value = "2021-09-14T10:16:43.906336Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:55.227339"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/700/?format=json"

# D2 ("1.18") must stay textual like the panel version string -- a leading
# apostrophe forces Excel to store it as text instead of coercing it to a
# floating point number.
$ws2.Range("D2").Value = "'1.18"

# --- Match header/index-column formatting used on the "data" sheet (bold,
# thin border, centered) by copying the existing style instead of rebuilding
# it from scratch so the workbook keeps reusing the same style record. ---
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Restore D2 to the default (unstyled) look -- the quote-prefix trick above
# nudges the cell into a temporary style; copying the plain format from an
# untouched cell clears that back out while leaving the text value in place.
$ws1.Range("A1").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

# --- Match page margins used throughout the workbook ---
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Keep "data" as the active/visible tab, as it was before the edit.
$ws1.Activate()
